$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is numeric-looking text (e.g. "228.93") need NumberFormat
# forced to Text ("@") before assignment, otherwise Excel auto-converts the
# literal string into a real number (matching native Excel typing semantics).
# The original (un-styled, style index 0) appearance is restored afterwards by
# copying the Style of an untouched plain cell (B2) back onto them, so the only
# net effect is the cell value/text -- no lingering style/format diff.
$plainStyle = $ws.Range("B2").Style
$forceTextCells = @("D5", "D6", "D8", "D9", "D10", "D13", "D14", "D15", "D19", "D20", "D22", "D25", "D26", "D27", "D28", "D29", "D30", "D32", "D33", "D38", "D39", "D41", "D43", "D44", "D45", "D46", "D48", "D51")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "37.362.78"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").Value = "2.046.18"
$ws.Range("E3").Value = "  -2.10%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "228.93"
$ws.Range("E5").Value = "  -1.99%  "
$ws.Range("D6").Value = "0.613"
$ws.Range("E6").Value = "  -2.03%  "
$ws.Range("D8").Value = "56.25"
$ws.Range("E8").Value = "  -3.72%  "
$ws.Range("D9").Value = "0.385"
$ws.Range("E9").Value = "  -2.55%  "
$ws.Range("D10").Value = "0.0817"
$ws.Range("E10").Value = "  +4.29%  "
$ws.Range("E11").Value = "  -1.97%  "
$ws.Range("D12").Value = "2.348.37"
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("D13").Value = "14.54"
$ws.Range("E13").Value = "  -3.82%  "
$ws.Range("D14").Value = "20.60"
$ws.Range("E14").Value = "  -3.05%  "
$ws.Range("D15").Value = "0.754"
$ws.Range("E15").Value = "  -3.25%  "
$ws.Range("E16").Value = "  -1.95%  "
$ws.Range("D17").Value = "2.038.18"
$ws.Range("E17").Value = "  -2.89%  "
$ws.Range("D18").Value = "37.244.99"
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("D19").Value = "6.06"
$ws.Range("E19").Value = "  -1.41%  "
$ws.Range("D20").Value = "69.74"
$ws.Range("E20").Value = "  -2.15%  "
$ws.Range("D21").Value = "0.0₃0861"
$ws.Range("E21").Value = "  +2.85%  "
$ws.Range("D22").Value = "225.96"
$ws.Range("E22").Value = "  -1.94%  "
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("D25").Value = "2.28"
$ws.Range("E25").Value = "  -4.83%  "
$ws.Range("D26").Value = "9.54"
$ws.Range("E26").Value = "  -2.51%  "
$ws.Range("D27").Value = "168.14"
$ws.Range("E27").Value = "  -1.89%  "
$ws.Range("D28").Value = "0.130"
$ws.Range("E28").Value = "  -4.33%  "
$ws.Range("D29").Value = "1.40"
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("D30").Value = "18.95"
$ws.Range("E30").Value = "  -2.97%  "
$ws.Range("E31").Value = "  -2.55%  "
$ws.Range("D32").Value = "4.52"
$ws.Range("E32").Value = "  -4.21%  "
$ws.Range("D33").Value = "0.0612"
$ws.Range("E33").Value = "  -3.41%  "
$ws.Range("E34").Value = "  -2.69%  "
$ws.Range("E35").Value = "  -5.18%  "
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("D38").Value = "3.18"
$ws.Range("E38").Value = "  -5.20%  "
$ws.Range("D39").Value = "5.39"
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("D40").Value = "1.506.59"
$ws.Range("E40").Value = "  +3.51%  "
$ws.Range("D41").Value = "0.0220"
$ws.Range("E41").Value = "  -6.59%  "
$ws.Range("E42").Value = "  -1.94%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "16.84"
$ws.Range("E43").Value = "  +0.90%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "96.12"
$ws.Range("E44").Value = "  -5.73%  "
$ws.Range("D45").Value = "0.0935"
$ws.Range("E45").Value = "  -4.13%  "
$ws.Range("D46").Value = "1.15"
$ws.Range("E46").Value = "  -2.83%  "
$ws.Range("E47").Value = "  -4.56%  "
$ws.Range("D48").Value = "7.15"
$ws.Range("E48").Value = "  -2.01%  "
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("D50").Value = "2.232.37"
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("D51").Value = "3.64"
$ws.Range("E51").Value = "  -11.85%  "

foreach ($addr in $forceTextCells) {
    $ws.Range($addr).Style = $plainStyle
}
